$wb = $excel.ActiveWorkbook

# Rename "pop" -> "population" and "pop_benelux" -> "population_benelux"
$wb.Worksheets.Item("pop").Name = "population"
$wb.Worksheets.Item("pop_benelux").Name = "population_benelux"
